# Update database values and change read_price algorithm:
# All the quarterly numeric figures in rows 11-27 (columns D:M) are reset.
# Rows that used to hold an all-zero placeholder ("هزینه کاهش ارزش دریافتنی‌ها"
# on row 15 and "سود (زیان) عملیات متوقف شده..." on row 23) now show a literal
# dash "-" instead of 0, while every other data row is zeroed out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dashRows = @(15, 23)
$dataRows = 11..27

foreach ($r in $dataRows) {
    if ($dashRows -contains $r) {
        $value = "-"
    } else {
        $value = 0
    }
    for ($col = 4; $col -le 13; $col++) {
        $ws.Cells.Item($r, $col).Value = $value
    }
}
